# Word COM-interop script
#
# 1) Several table-cell placeholders ("{entry_time}", "{entry_gate}",
#    "{exit_time}", "{exit_gate}", "{paid_amt}", "{offender_name}",
#    "{offender_contact}", "{offender_vehicle}", "{offender_store}")
#    were previously split across three runs (the spell-checker had
#    wrapped the middle portion with proofErr spellStart/spellEnd).
#    A Find/Replace over the whole visible string collapses those three
#    runs (and drops the now-orphaned proofErr markers) back into a
#    single run, matching how Word normally stores an untouched field.
#
# 2) The "{images}" placeholder (a single run) is changed into a
#    mail-merge-style "{%images%}" token that must be stored as five
#    separate runs: "{", "%", "images", "%", "}". Because the runs all
#    share identical (empty) formatting, a plain Find/Replace or
#    Range.InsertAfter would just be re-coalesced into one run on save,
#    so we instead replace the found range's contents with literal
#    OOXML (Range.InsertXML) that explicitly encodes five runs inside
#    the original paragraph (keeping the paragraph's own attributes
#    unchanged).

$d = $word.ActiveDocument

# --- 1) Merge the split placeholder runs back into single runs -------------
$fields = @(
    "entry_time",
    "entry_gate",
    "exit_time",
    "exit_gate",
    "paid_amt",
    "offender_name",
    "offender_contact",
    "offender_vehicle",
    "offender_store"
)

foreach ($f in $fields) {
    $token = "{" + $f + "}"
    $d.Content.Find.Execute($token, $true, $false, $false, $false, $false, $true, 1, $false, $token, 2)
}

# --- 2) Split "{images}" into "{", "%", "images", "%", "}" runs ------------
# A plain Find/Replace (or Range.InsertAfter) would leave the five pieces
# merged back into a single run on save because they all share identical
# (empty) run formatting, so the paragraph's own run content is instead
# replaced with literal OOXML describing five discrete runs. The
# paragraph's own attributes (paraId/textId/rsid*) are carried over
# unchanged, exactly as in the original document.
$rng = $d.Content
$found = $rng.Find.Execute("{images}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($rng.Start, $rng.End)

    $xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
           'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
           'w14:paraId="4A2CDFC4" w14:textId="19D18148" w:rsidR="00052777" ' +
           'w:rsidRDefault="001373FD" w:rsidP="3B00479A">' +
           '<w:r><w:t>{</w:t></w:r>' +
           '<w:r><w:t>%</w:t></w:r>' +
           '<w:r><w:t>images</w:t></w:r>' +
           '<w:r><w:t>%</w:t></w:r>' +
           '<w:r><w:t>}</w:t></w:r>' +
           '</w:p></pkg:xmlData>'

    $target.InsertXML($xml)
}
